$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that is bumped by one day
# for every data row (rows 2-86) as part of the automatic update.
for ($row = 2; $row -le 86; $row++) {
    $ws.Cells.Item($row, 3).Value = 46062
}
